$wb = $excel.ActiveWorkbook

# Sheet tab names whose A1 cell header needs "id='...'" changed to "class='...'"
$sheetNames = @(
    "!!Model",
    "!!Taxon",
    "!!Environment",
    "!!Submodels",
    "!!Compartments",
    "!!Species types",
    "!!Species",
    "!!Initial species concentrations",
    "!!Observables",
    "!!Functions",
    "!!Reactions",
    "!!Rate laws",
    "!!dFBA objectives",
    "!!dFBA objective reactions",
    "!!dFBA objective species",
    "!!Parameters",
    "!!Stop conditions",
    "!!Observations",
    "!!Observation sets",
    "!!Conclusions",
    "!!References",
    "!!Authors",
    "!!Changes"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $cell = $ws.Range("A1")
    $oldValue = $cell.Value
    if ($oldValue -like "*ObjTables type='Data' id='*") {
        $newValue = $oldValue -replace "type='Data' id='", "type='Data' class='"
        $cell.Value = $newValue
    }
}
